# Recomputed transition-probability matrix (team_specific_matrix/Lipscomb_B)
# after simulating more games / updated simulate-game logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2455357142857143
$ws.Range("C2").Value = 0.4910714285714285
$ws.Range("J2").Value = 0.004464285714285714
$ws.Range("P2").Value = 0.1741071428571428
$ws.Range("S2").Value = 0.08482142857142858
$ws.Range("C3").Value = 0.02654867256637168
$ws.Range("J3").Value = 0.01769911504424779
$ws.Range("P3").Value = 0.7787610619469026
$ws.Range("S3").Value = 0.1769911504424779
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7560975609756098
$ws.Range("S4").Value = 0.2195121951219512
$ws.Range("B6").Value = 0.04678362573099415
$ws.Range("D6").Value = 0.02339181286549707
$ws.Range("F6").Value = 0.03508771929824561
$ws.Range("J6").Value = 0.2865497076023392
$ws.Range("O6").Value = 0.01754385964912281
$ws.Range("Q6").Value = 0.1754385964912281
$ws.Range("R6").Value = 0.07017543859649122
$ws.Range("S6").Value = 0.3450292397660819
$ws.Range("B7").Value = 0.08247422680412371
$ws.Range("D7").Value = 0.02061855670103093
$ws.Range("F7").Value = 0.04123711340206185
$ws.Range("J7").Value = 0.1237113402061856
$ws.Range("O7").Value = 0.01030927835051546
$ws.Range("Q7").Value = 0.2216494845360825
$ws.Range("R7").Value = 0.08247422680412371
$ws.Range("S7").Value = 0.4175257731958763
$ws.Range("B8").Value = 0.06693711967545639
$ws.Range("D8").Value = 0.01825557809330629
$ws.Range("E8").Value = 0.002028397565922921
$ws.Range("F8").Value = 0.02636916835699797
$ws.Range("J8").Value = 0.1318458417849899
$ws.Range("O8").Value = 0.004056795131845842
$ws.Range("Q8").Value = 0.1906693711967546
$ws.Range("R8").Value = 0.1176470588235294
$ws.Range("S8").Value = 0.4421906693711968
$ws.Range("B9").Value = 0.06382978723404255
$ws.Range("D9").Value = 0.02553191489361702
$ws.Range("E9").Value = 0.00425531914893617
$ws.Range("F9").Value = 0.05957446808510639
$ws.Range("J9").Value = 0.1276595744680851
$ws.Range("O9").Value = 0.02127659574468085
$ws.Range("Q9").Value = 0.2553191489361702
$ws.Range("R9").Value = 0.08936170212765958
$ws.Range("S9").Value = 0.3531914893617021
$ws.Range("B10").Value = 0.07963118189438391
$ws.Range("D10").Value = 0.01676445934618609
$ws.Range("E10").Value = 0.0008382229673093043
$ws.Range("F10").Value = 0.06286672254819782
$ws.Range("J10").Value = 0.1089689857502096
$ws.Range("O10").Value = 0.01508801341156748
$ws.Range("Q10").Value = 0.2388935456831517
$ws.Range("R10").Value = 0.09639564124056998
$ws.Range("S10").Value = 0.3805532271584242
$ws.Range("G11").Value = 0.150197628458498
$ws.Range("J11").Value = 0.07114624505928854
$ws.Range("K11").Value = 0.2015810276679842
$ws.Range("L11").Value = 0.5691699604743083
$ws.Range("S11").Value = 0.007905138339920948
$ws.Range("G12").Value = 0.815068493150685
$ws.Range("J12").Value = 0.1506849315068493
$ws.Range("L12").Value = 0.02054794520547945
$ws.Range("S12").Value = 0.0136986301369863
$ws.Range("G13").Value = 0.82
$ws.Range("J13").Value = 0.18
$ws.Range("F15").Value = 0.005128205128205128
$ws.Range("H15").Value = 0.1948717948717949
$ws.Range("I15").Value = 0.09743589743589744
$ws.Range("J15").Value = 0.358974358974359
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.01025641025641026
$ws.Range("O15").Value = 0.04102564102564103
$ws.Range("S15").Value = 0.2256410256410256
$ws.Range("F16").Value = 0.0130718954248366
$ws.Range("H16").Value = 0.1895424836601307
$ws.Range("I16").Value = 0.1241830065359477
$ws.Range("J16").Value = 0.4313725490196079
$ws.Range("K16").Value = 0.09803921568627451
$ws.Range("M16").Value = 0.03267973856209151
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.05228758169934641
$ws.Range("F17").Value = 0.01577909270216963
$ws.Range("H17").Value = 0.1715976331360947
$ws.Range("I17").Value = 0.1400394477317554
$ws.Range("J17").Value = 0.398422090729783
$ws.Range("K17").Value = 0.09664694280078895
$ws.Range("M17").Value = 0.01972386587771203
$ws.Range("O17").Value = 0.05325443786982249
$ws.Range("S17").Value = 0.1045364891518738
$ws.Range("F18").Value = 0.004587155963302753
$ws.Range("H18").Value = 0.2201834862385321
$ws.Range("I18").Value = 0.07339449541284404
$ws.Range("J18").Value = 0.426605504587156
$ws.Range("K18").Value = 0.1146788990825688
$ws.Range("M18").Value = 0.02293577981651376
$ws.Range("O18").Value = 0.05045871559633028
$ws.Range("S18").Value = 0.0871559633027523
$ws.Range("F19").Value = 0.01198630136986301
$ws.Range("H19").Value = 0.25
$ws.Range("I19").Value = 0.09674657534246575
$ws.Range("J19").Value = 0.3647260273972603
$ws.Range("K19").Value = 0.08304794520547945
$ws.Range("M19").Value = 0.02397260273972603
$ws.Range("O19").Value = 0.0761986301369863
$ws.Range("S19").Value = 0.09332191780821918
